$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values (recalculated statistics for updated year ranges) ---
$ws.Range("G4").Value = 0.0133981016431697
$ws.Range("L4").Value = 0.00181
$ws.Range("G5").Value = 0.0133981016431697
$ws.Range("L5").Value = 0.00181
$ws.Range("G6").Value = 0.672624043048015
$ws.Range("L6").Value = 0.04225
$ws.Range("G7").Value = 0.672624043048015
$ws.Range("L7").Value = 0.04225
$ws.Range("G8").Value = 0.725547730102311
$ws.Range("G9").Value = 0.725547730102311
$ws.Range("G20").Value = 0.0182368548944615
$ws.Range("L20").Value = 0.00485
$ws.Range("G21").Value = 0.0182368548944615
$ws.Range("L21").Value = 0.00485
$ws.Range("G22").Value = 0.702254571293116
$ws.Range("G23").Value = 0.702254571293116

# --- Append new rows 30-45 (2019 - 2023 period results) ---
$ws.Range("A30").Value = 'Koitiata at Beamish Rd'
$ws.Range("B30").Value = 'DRP (95th Percentile)'
$ws.Range("C30").Value = 'D'
$ws.Range("D30").Value = '2019 - 2023'
$ws.Range("E30").Value = 'RepSite'
$ws.Range("F30").Value = 0.021
$ws.Range("G30").Value = 0.0273962264150943
$ws.Range("H30").Value = 0.127
$ws.Range("I30").Value = 0.11925
$ws.Range("L30").Value = 0.021
$ws.Range("M30").Value = 0.02449
$ws.Range("N30").Value = 0.0551
$ws.Range("O30").Value = 1785275
$ws.Range("P30").Value = 5558196
$ws.Range("Q30").Value = 'Rangitikei District'
$ws.Range("R30").Value = 'Rangitīkei-Turakina'
$ws.Range("S30").Value = 'Southern Whanganui Lakes'
$ws.Range("T30").Value = 'West_5'
$ws.Range("U30").Value = 'mg/L'
$ws.Range("A31").Value = 'Koitiata at Beamish Rd'
$ws.Range("B31").Value = 'DRP (Median)'
$ws.Range("C31").Value = 'D'
$ws.Range("D31").Value = '2019 - 2023'
$ws.Range("E31").Value = 'RepSite'
$ws.Range("F31").Value = 0.021
$ws.Range("G31").Value = 0.0273962264150943
$ws.Range("H31").Value = 0.127
$ws.Range("I31").Value = 0.11925
$ws.Range("L31").Value = 0.021
$ws.Range("M31").Value = 0.02449
$ws.Range("N31").Value = 0.0551
$ws.Range("O31").Value = 1785275
$ws.Range("P31").Value = 5558196
$ws.Range("Q31").Value = 'Rangitikei District'
$ws.Range("R31").Value = 'Rangitīkei-Turakina'
$ws.Range("S31").Value = 'Southern Whanganui Lakes'
$ws.Range("T31").Value = 'West_5'
$ws.Range("U31").Value = 'mg/L'
$ws.Range("A32").Value = 'Koitiata at Beamish Rd'
$ws.Range("B32").Value = 'E coli (>260)'
$ws.Range("C32").Value = 'B'
$ws.Range("D32").Value = '2019 - 2023'
$ws.Range("E32").Value = 'RepSite'
$ws.Range("F32").Value = 120
$ws.Range("G32").Value = 1048.60228962081
$ws.Range("H32").Value = 19935.0336489237
$ws.Range("I32").Value = 4112.25
$ws.Range("J32").Value = 20.7547169811321
$ws.Range("K32").Value = 24.5283018867925
$ws.Range("L32").Value = 94
$ws.Range("M32").Value = 677.27
$ws.Range("N32").Value = 2994.78
$ws.Range("O32").Value = 1785275
$ws.Range("P32").Value = 5558196
$ws.Range("Q32").Value = 'Rangitikei District'
$ws.Range("R32").Value = 'Rangitīkei-Turakina'
$ws.Range("S32").Value = 'Southern Whanganui Lakes'
$ws.Range("T32").Value = 'West_5'
$ws.Range("U32").Value = '% exceedances over 260/100 mL'
$ws.Range("A33").Value = 'Koitiata at Beamish Rd'
$ws.Range("B33").Value = 'E coli (>540)'
$ws.Range("C33").Value = 'D'
$ws.Range("D33").Value = '2019 - 2023'
$ws.Range("E33").Value = 'RepSite'
$ws.Range("F33").Value = 120
$ws.Range("G33").Value = 1048.60228962081
$ws.Range("H33").Value = 19935.0336489237
$ws.Range("I33").Value = 4112.25
$ws.Range("J33").Value = 20.7547169811321
$ws.Range("K33").Value = 24.5283018867925
$ws.Range("L33").Value = 94
$ws.Range("M33").Value = 677.27
$ws.Range("N33").Value = 2994.78
$ws.Range("O33").Value = 1785275
$ws.Range("P33").Value = 5558196
$ws.Range("Q33").Value = 'Rangitikei District'
$ws.Range("R33").Value = 'Rangitīkei-Turakina'
$ws.Range("S33").Value = 'Southern Whanganui Lakes'
$ws.Range("T33").Value = 'West_5'
$ws.Range("U33").Value = '% exceedances over 540/100 mL'
$ws.Range("A34").Value = 'Koitiata at Beamish Rd'
$ws.Range("B34").Value = 'E coli (Median)'
$ws.Range("C34").Value = 'A'
$ws.Range("D34").Value = '2019 - 2023'
$ws.Range("E34").Value = 'RepSite'
$ws.Range("F34").Value = 120
$ws.Range("G34").Value = 1048.60228962081
$ws.Range("H34").Value = 19935.0336489237
$ws.Range("I34").Value = 4112.25
$ws.Range("J34").Value = 20.7547169811321
$ws.Range("K34").Value = 24.5283018867925
$ws.Range("L34").Value = 94
$ws.Range("M34").Value = 677.27
$ws.Range("N34").Value = 2994.78
$ws.Range("O34").Value = 1785275
$ws.Range("P34").Value = 5558196
$ws.Range("Q34").Value = 'Rangitikei District'
$ws.Range("R34").Value = 'Rangitīkei-Turakina'
$ws.Range("S34").Value = 'Southern Whanganui Lakes'
$ws.Range("T34").Value = 'West_5'
$ws.Range("U34").Value = 'E. coli/100 mL'
$ws.Range("A35").Value = 'Koitiata at Beamish Rd'
$ws.Range("B35").Value = 'E coli (95th Percentile)'
$ws.Range("C35").Value = 'E'
$ws.Range("D35").Value = '2019 - 2023'
$ws.Range("E35").Value = 'RepSite'
$ws.Range("F35").Value = 120
$ws.Range("G35").Value = 1048.60228962081
$ws.Range("H35").Value = 19935.0336489237
$ws.Range("I35").Value = 4112.25
$ws.Range("J35").Value = 20.7547169811321
$ws.Range("K35").Value = 24.5283018867925
$ws.Range("L35").Value = 94
$ws.Range("M35").Value = 677.27
$ws.Range("N35").Value = 2994.78
$ws.Range("O35").Value = 1785275
$ws.Range("P35").Value = 5558196
$ws.Range("Q35").Value = 'Rangitikei District'
$ws.Range("R35").Value = 'Rangitīkei-Turakina'
$ws.Range("S35").Value = 'Southern Whanganui Lakes'
$ws.Range("T35").Value = 'West_5'
$ws.Range("U35").Value = 'E. coli/100 mL'
$ws.Range("A36").Value = 'Koitiata at Beamish Rd'
$ws.Range("B36").Value = 'Ammoniacal-N (95th Percentile)'
$ws.Range("C36").Value = 'B'
$ws.Range("D36").Value = '2019 - 2023'
$ws.Range("E36").Value = 'RepSite'
$ws.Range("F36").Value = 0.01137
$ws.Range("G36").Value = 0.0736703165436941
$ws.Range("H36").Value = 1.7508335155888
$ws.Range("I36").Value = 0.36227
$ws.Range("L36").Value = 0.00685
$ws.Range("M36").Value = 0.04769
$ws.Range("N36").Value = 0.14058
$ws.Range("O36").Value = 1785275
$ws.Range("P36").Value = 5558196
$ws.Range("Q36").Value = 'Rangitikei District'
$ws.Range("R36").Value = 'Rangitīkei-Turakina'
$ws.Range("S36").Value = 'Southern Whanganui Lakes'
$ws.Range("T36").Value = 'West_5'
$ws.Range("U36").Value = 'mg NH4-N/L'
$ws.Range("A37").Value = 'Koitiata at Beamish Rd'
$ws.Range("B37").Value = 'Ammoniacal-N (Median)'
$ws.Range("C37").Value = 'A'
$ws.Range("D37").Value = '2019 - 2023'
$ws.Range("E37").Value = 'RepSite'
$ws.Range("F37").Value = 0.01137
$ws.Range("G37").Value = 0.0736703165436941
$ws.Range("H37").Value = 1.7508335155888
$ws.Range("I37").Value = 0.36227
$ws.Range("L37").Value = 0.00685
$ws.Range("M37").Value = 0.04769
$ws.Range("N37").Value = 0.14058
$ws.Range("O37").Value = 1785275
$ws.Range("P37").Value = 5558196
$ws.Range("Q37").Value = 'Rangitikei District'
$ws.Range("R37").Value = 'Rangitīkei-Turakina'
$ws.Range("S37").Value = 'Southern Whanganui Lakes'
$ws.Range("T37").Value = 'West_5'
$ws.Range("U37").Value = 'mg NH4-N/L'
$ws.Range("A38").Value = 'Koitiata at Beamish Rd'
$ws.Range("B38").Value = 'Nitrate-N (95th Percentile)'
$ws.Range("C38").Value = 'B'
$ws.Range("D38").Value = '2019 - 2023'
$ws.Range("E38").Value = 'RepSite'
$ws.Range("F38").Value = 0.623
$ws.Range("G38").Value = 0.784254571293116
$ws.Range("H38").Value = 2.28
$ws.Range("I38").Value = 1.9515
$ws.Range("L38").Value = 0.1255
$ws.Range("M38").Value = 1.3798
$ws.Range("N38").Value = 1.7826
$ws.Range("O38").Value = 1785275
$ws.Range("P38").Value = 5558196
$ws.Range("Q38").Value = 'Rangitikei District'
$ws.Range("R38").Value = 'Rangitīkei-Turakina'
$ws.Range("S38").Value = 'Southern Whanganui Lakes'
$ws.Range("T38").Value = 'West_5'
$ws.Range("U38").Value = 'mg NO3-N/L'
$ws.Range("A39").Value = 'Koitiata at Beamish Rd'
$ws.Range("B39").Value = 'Nitrate-N (Median)'
$ws.Range("C39").Value = 'A'
$ws.Range("D39").Value = '2019 - 2023'
$ws.Range("E39").Value = 'RepSite'
$ws.Range("F39").Value = 0.623
$ws.Range("G39").Value = 0.784254571293116
$ws.Range("H39").Value = 2.28
$ws.Range("I39").Value = 1.9515
$ws.Range("L39").Value = 0.1255
$ws.Range("M39").Value = 1.3798
$ws.Range("N39").Value = 1.7826
$ws.Range("O39").Value = 1785275
$ws.Range("P39").Value = 5558196
$ws.Range("Q39").Value = 'Rangitikei District'
$ws.Range("R39").Value = 'Rangitīkei-Turakina'
$ws.Range("S39").Value = 'Southern Whanganui Lakes'
$ws.Range("T39").Value = 'West_5'
$ws.Range("U39").Value = 'mg NO3-N/L'
$ws.Range("A40").Value = 'Koitiata at Beamish Rd'
$ws.Range("B40").Value = 'Soluble Inorganic Nitrogen (95th Percentile)'
$ws.Range("D40").Value = '2019 - 2023'
$ws.Range("E40").Value = 'RepSite'
$ws.Range("F40").Value = 0.73
$ws.Range("G40").Value = 0.939433962264151
$ws.Range("H40").Value = 4.75
$ws.Range("I40").Value = 2.348
$ws.Range("L40").Value = 0.13
$ws.Range("M40").Value = 1.7049
$ws.Range("N40").Value = 2.0508
$ws.Range("O40").Value = 1785275
$ws.Range("P40").Value = 5558196
$ws.Range("Q40").Value = 'Rangitikei District'
$ws.Range("R40").Value = 'Rangitīkei-Turakina'
$ws.Range("S40").Value = 'Southern Whanganui Lakes'
$ws.Range("T40").Value = 'West_5'
$ws.Range("U40").Value = 'g/m3'
$ws.Range("A41").Value = 'Koitiata at Beamish Rd'
$ws.Range("B41").Value = 'Soluble Inorganic Nitrogen (Median)'
$ws.Range("D41").Value = '2019 - 2023'
$ws.Range("E41").Value = 'RepSite'
$ws.Range("F41").Value = 0.73
$ws.Range("G41").Value = 0.939433962264151
$ws.Range("H41").Value = 4.75
$ws.Range("I41").Value = 2.348
$ws.Range("L41").Value = 0.13
$ws.Range("M41").Value = 1.7049
$ws.Range("N41").Value = 2.0508
$ws.Range("O41").Value = 1785275
$ws.Range("P41").Value = 5558196
$ws.Range("Q41").Value = 'Rangitikei District'
$ws.Range("R41").Value = 'Rangitīkei-Turakina'
$ws.Range("S41").Value = 'Southern Whanganui Lakes'
$ws.Range("T41").Value = 'West_5'
$ws.Range("U41").Value = 'g/m3'
$ws.Range("A42").Value = 'Koitiata at Beamish Rd'
$ws.Range("B42").Value = 'Total Nitrogen (95th Percentile)'
$ws.Range("D42").Value = '2019 - 2023'
$ws.Range("E42").Value = 'RepSite'
$ws.Range("F42").Value = 1.64
$ws.Range("G42").Value = 1.81377358490566
$ws.Range("H42").Value = 6.97
$ws.Range("I42").Value = 3.61
$ws.Range("L42").Value = 0.93
$ws.Range("M42").Value = 2.4794
$ws.Range("N42").Value = 3.159
$ws.Range("O42").Value = 1785275
$ws.Range("P42").Value = 5558196
$ws.Range("Q42").Value = 'Rangitikei District'
$ws.Range("R42").Value = 'Rangitīkei-Turakina'
$ws.Range("S42").Value = 'Southern Whanganui Lakes'
$ws.Range("T42").Value = 'West_5'
$ws.Range("U42").Value = 'g/m3'
$ws.Range("A43").Value = 'Koitiata at Beamish Rd'
$ws.Range("B43").Value = 'Total Nitrogen (Median)'
$ws.Range("D43").Value = '2019 - 2023'
$ws.Range("E43").Value = 'RepSite'
$ws.Range("F43").Value = 1.64
$ws.Range("G43").Value = 1.81377358490566
$ws.Range("H43").Value = 6.97
$ws.Range("I43").Value = 3.61
$ws.Range("L43").Value = 0.93
$ws.Range("M43").Value = 2.4794
$ws.Range("N43").Value = 3.159
$ws.Range("O43").Value = 1785275
$ws.Range("P43").Value = 5558196
$ws.Range("Q43").Value = 'Rangitikei District'
$ws.Range("R43").Value = 'Rangitīkei-Turakina'
$ws.Range("S43").Value = 'Southern Whanganui Lakes'
$ws.Range("T43").Value = 'West_5'
$ws.Range("U43").Value = 'g/m3'
$ws.Range("A44").Value = 'Koitiata at Beamish Rd'
$ws.Range("B44").Value = 'Total Phosphorus (95th Percentile)'
$ws.Range("D44").Value = '2019 - 2023'
$ws.Range("E44").Value = 'RepSite'
$ws.Range("F44").Value = 0.04
$ws.Range("G44").Value = 0.32
$ws.Range("H44").Value = 12.2
$ws.Range("I44").Value = 0.5016
$ws.Range("L44").Value = 0.032
$ws.Range("M44").Value = 0.13435
$ws.Range("N44").Value = 0.40884
$ws.Range("O44").Value = 1785275
$ws.Range("P44").Value = 5558196
$ws.Range("Q44").Value = 'Rangitikei District'
$ws.Range("R44").Value = 'Rangitīkei-Turakina'
$ws.Range("S44").Value = 'Southern Whanganui Lakes'
$ws.Range("T44").Value = 'West_5'
$ws.Range("U44").Value = 'g/m3'
$ws.Range("A45").Value = 'Koitiata at Beamish Rd'
$ws.Range("B45").Value = 'Total Phosphorus (Median)'
$ws.Range("D45").Value = '2019 - 2023'
$ws.Range("E45").Value = 'RepSite'
$ws.Range("F45").Value = 0.04
$ws.Range("G45").Value = 0.32
$ws.Range("H45").Value = 12.2
$ws.Range("I45").Value = 0.5016
$ws.Range("L45").Value = 0.032
$ws.Range("M45").Value = 0.13435
$ws.Range("N45").Value = 0.40884
$ws.Range("O45").Value = 1785275
$ws.Range("P45").Value = 5558196
$ws.Range("Q45").Value = 'Rangitikei District'
$ws.Range("R45").Value = 'Rangitīkei-Turakina'
$ws.Range("S45").Value = 'Southern Whanganui Lakes'
$ws.Range("T45").Value = 'West_5'
$ws.Range("U45").Value = 'g/m3'

